# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型
# sheets to match the latest scrape (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9867
$ws1.Range("F8").Value = 478
$ws1.Range("F9").Value = 722
$ws1.Range("F13").Value = 3051
$ws1.Range("F14").Value = 2306
$ws1.Range("F16").Value = 2002
$ws1.Range("F20").Value = 1571
$ws1.Range("F21").Value = 429
$ws1.Range("F22").Value = 40
$ws1.Range("F23").Value = 205
$ws1.Range("F24").Value = 227
$ws1.Range("F26").Value = 355
$ws1.Range("F28").Value = 336
$ws1.Range("F29").Value = 543
$ws1.Range("F31").Value = 188
$ws1.Range("F32").Value = 1556
$ws1.Range("F33").Value = 249
$ws1.Range("F34").Value = 1579
$ws1.Range("F35").Value = 77
$ws1.Range("F36").Value = 383
$ws1.Range("F37").Value = 38
$ws1.Range("F38").Value = 412
$ws1.Range("F39").Value = 853
$ws1.Range("F41").Value = 332

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9867
$ws4.Range("F10").Value = 478
$ws4.Range("F11").Value = 722
$ws4.Range("F15").Value = 3051
$ws4.Range("F16").Value = 2306
$ws4.Range("F18").Value = 2002
$ws4.Range("F22").Value = 1571
$ws4.Range("F23").Value = 429
$ws4.Range("F24").Value = 40
$ws4.Range("F25").Value = 205
$ws4.Range("F26").Value = 227
$ws4.Range("F28").Value = 355
$ws4.Range("F30").Value = 336
$ws4.Range("F31").Value = 543
$ws4.Range("F36").Value = 188
$ws4.Range("F37").Value = 1556
$ws4.Range("F39").Value = 249
$ws4.Range("F40").Value = 1579
$ws4.Range("F41").Value = 77
$ws4.Range("F43").Value = 383
$ws4.Range("F44").Value = 38
$ws4.Range("F45").Value = 412
$ws4.Range("F46").Value = 854
$ws4.Range("F48").Value = 332
